# Update database and change read_price algorithm
# - Drop the oldest quarter column (D, "فصل دوم منتهی به 1399/06") by deleting
#   the whole column, which shifts every later quarter one column to the left.
# - Amend the publish-date note for "فصل چهارم منتهی به 1400/12"
#   (now column I after the shift) to reflect a later republish date.
# - Append the new quarter "فصل چهارم منتهی به 1401/12" in the freed-up
#   column M, together with its figures and publish date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove oldest quarter column; everything to the right shifts left.
$ws.Columns("D").Delete()

# 2. The quarter that lands on column I ("فصل چهارم منتهی به 1400/12") was
#    republished later, so update its publish-date label in place.
$ws.Range("I9").Value = "1402-02-28 (7)"

# 3. Fill in the new quarter in column M ("فصل چهارم منتهی به 1401/12").
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-28"

$ws.Range("M11").Value = 5682138
$ws.Range("M12").Value = -4616217
$ws.Range("M13").Value = 1065921
$ws.Range("M14").Value = -237125
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 131170
$ws.Range("M17").Value = 959966
$ws.Range("M18").Value = -180950
$ws.Range("M19").Value = -64874
$ws.Range("M20").Value = 714142
$ws.Range("M21").Value = -36093
$ws.Range("M22").Value = 678049
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 678049
$ws.Range("M25").Value = 52
$ws.Range("M26").Value = 13000000
$ws.Range("M27").Value = 52

# 4. Keep the "wide" column styling pattern (every third data column is
#    wider) now that column M is the newest, right-most quarter.
$ws.Columns("M").ColumnWidth = 31
